# "Penalty Reward System" attempt: drop the stale mid-history weeks/months
# that were superseded once the 2024 data arrived, so each sheet's trend
# starts right where the older, now-irrelevant rows used to sit.
$wb = $excel.ActiveWorkbook

# Weekly Quantity: remove the 7 weeks from 2023-06-11 through 2023-07-30
# (rows 16-22), shifting the 2024 rows (formerly 23-30) up to 16-23.
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("A16:B22").EntireRow.Delete()

# Monthly Trend: remove the 2 months 2023-07 and 2023-08 (rows 6-7),
# shifting the 2024 rows (formerly 8-10) up to 6-8.
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("A6:B7").EntireRow.Delete()
